# Update 2021 conference championship game stats (simulated + logged)
# for the Kansas City Chiefs "Players Data" workbook.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")
$rushing.Range("C2").Value = 20   # P.Mahomes 1DATT
$rushing.Range("D2").Value = 21   # P.Mahomes 2DATT
$rushing.Range("F2").Value = 17   # P.Mahomes RZATT

$rushing.Range("C4").Value = 82   # C.Edwards-Helaire 1DATT
$rushing.Range("D4").Value = 48   # C.Edwards-Helaire 2DATT

$rushing.Range("C6").Value = 27   # J.McKinnon 1DATT
$rushing.Range("D6").Value = 13   # J.McKinnon 2DATT
$rushing.Range("F6").Value = 8    # J.McKinnon RZATT

$rushing.Range("E8").Value = 7    # M.Burton 3DATT

$rushing.Range("C10").Value = 11  # M.Hardman 1DATT
$rushing.Range("D10").Value = 1   # M.Hardman 2DATT

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")
$receiving.Range("C2").Value = 41   # C.Edwards-Helaire Short Target
$receiving.Range("D2").Value = 31   # C.Edwards-Helaire Short Comp

$receiving.Range("C4").Value = 19   # J.McKinnon Short Target
$receiving.Range("D4").Value = 15   # J.McKinnon Short Comp
$receiving.Range("G4").Value = 6    # J.McKinnon RZ Target
$receiving.Range("H4").Value = 5    # J.McKinnon RZ Comp

$receiving.Range("C7").Value = 142  # T.Hill Short Target
$receiving.Range("D7").Value = 113  # T.Hill Short Comp
$receiving.Range("E7").Value = 42   # T.Hill Deep Target
$receiving.Range("F7").Value = 18   # T.Hill Deep Comp
$receiving.Range("G7").Value = 32   # T.Hill RZ Target
$receiving.Range("H7").Value = 22   # T.Hill RZ Comp

$receiving.Range("C8").Value = 67   # M.Hardman Short Target
$receiving.Range("D8").Value = 52   # M.Hardman Short Comp
$receiving.Range("E8").Value = 20   # M.Hardman Deep Target
$receiving.Range("F8").Value = 11   # M.Hardman Deep Comp
$receiving.Range("G8").Value = 15   # M.Hardman RZ Target
$receiving.Range("H8").Value = 8    # M.Hardman RZ Comp

$receiving.Range("C9").Value = 53   # B.Pringle Short Target
$receiving.Range("D9").Value = 37   # B.Pringle Short Comp
$receiving.Range("E9").Value = 16   # B.Pringle Deep Target

$receiving.Range("C10").Value = 35  # D.Robinson Short Target

$receiving.Range("C13").Value = 136 # T.Kelce Short Target
$receiving.Range("D13").Value = 99  # T.Kelce Short Comp
$receiving.Range("G13").Value = 24  # T.Kelce RZ Target
$receiving.Range("H13").Value = 20  # T.Kelce RZ Comp
